# Applies the three genuine text-content changes from the commit:
#  1. "I’m Chief Scientific Officer ..." -> "I am Chief Scientific Officer ..."
#  2. Insert a new sentence about R&D projects into the "ten years of experience" paragraph.
#  3. Rewrite the "I’m a truly passionate ..." paragraph.
#
# (The rest of the underlying diff is Word's automatic proofing markup -
#  <w:proofErr> spell/grammar-check wrappers and run re-splitting around
#  them - which carries no visible text change, so it is not reproduced
#  here.)

$d = $word.ActiveDocument

# 1) "I'm" -> "I am" for the Chief Scientific Officer sentence.
$d.Content.Find.Execute(
    "I’m Chief Scientific Officer at Symanto",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "I am Chief Scientific Officer at Symanto",
    2
) | Out-Null

# 2) Insert new sentence about R&D projects before "Those projects allowed ...".
$d.Content.Find.Execute(
    "information retrieval. Those projects allowed me to acquire",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "information retrieval. I have participated in more than ten R&D projects with national and international funding, in some of which I have been head and responsible of research. Those projects allowed me to acquire",
    2
) | Out-Null

# 3) Rewrite the "truly passionate" paragraph.
$d.Content.Find.Execute(
    "I’m a truly passionate, focused, and creative researcher. I like to read blogs and scientific papers about computational linguistics, data mining, and deep learning to keep up with the latest advances. While I enjoy all aspects of my job, I think my favorite stage of a project is defining the main data pipeline to solve a concrete objective. I also find very stimulating the part of the result analysis and pipeline tuning.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "I am a passionate, focused, and creative researcher. I like to read blogs and scientific papers about computational linguistics, data mining, and deep learning. In addition, I contribute to the scientific community by publishing articles in prestigious conferences and journals, and organizing shared tasks in national and international conferences.",
    2
) | Out-Null
